# Add the "Rock your Life" mentoring program row (Resnjanskij et al. 2021)
# to the programs worksheet, following the pattern of the existing rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rightQuote = [char]0x2019

$description = "Rock your Life is a mentoring program, where university students take on mentees from the lowest-track secondary schools (Hauptschule). The program aims at providing career guidance, establishing career visions, and fostering self-esteem and trust in the mentees" + $rightQuote + " own skills and abilities"

# New row (51): Links, Papers, program_identifier, program_name, year,
# category, average_age_beneficiary, short_description
$ws.Range("J51").Value = "https://ideas.repec.org/p/ces/ceswps/_8870.html"
$ws.Hyperlinks.Add($ws.Range("J51"), "https://ideas.repec.org/p/ces/ceswps/_8870.html") | Out-Null
$ws.Range("J51").Style = "Link"

$ws.Range("I51").Value = "Resnjanskij et al. (2021)"

$ws.Range("A51").Value = "rockYourLife"
$ws.Range("B51").Value = "Mentoring Program Rock your Life"

$ws.Range("C51").Value = 2017
$ws.Range("D51").Value = "Education"
$ws.Range("E51").Value = 14

$ws.Range("F51").Value = $description
$ws.Range("F51").WrapText = $true

# Match row height used for similarly long wrapped descriptions.
$ws.Rows.Item(51).RowHeight = 105

# Update the view so the new row is visible / selected, mirroring the diff.
$ws.Range("F51").Select() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollColumn = 2
$win.ScrollRow = 48
